# Commit: "Tue, Apr 07, 2020  5:06:14 AM"
#
# 1) Swap the table style used by the table on slide 16 to a different
#    built-in PowerPoint table style.
# 2) Re-colour the deck's theme ("Integral") so its 10 theme colours
#    (everything except dk1/lt1, which are already black/white in both
#    themes) match the stock "Office Theme" palette.

$p = $ppt.ActivePresentation

# --- 1. Table style -------------------------------------------------
$slide = $p.Slides.Item(16)
$tableShape = $slide.Shapes.Item(3)
$table = $tableShape.Table
$table.ApplyStyle("{ABD42441-0A3E-40BE-B0B5-30C05C76F374}")

# --- 2. Theme colours -------------------------------------------------
# Item indices on ThemeColorScheme: 1 dk1, 2 lt1, 3 dk2, 4 lt2,
# 5-10 accent1-6, 11 hlink, 12 folHlink.
function Set-ThemeColor($scheme, $index, $r, $g, $b) {
    $scheme.Item($index).RGB = $r + ($g * 256) + ($b * 65536)
}

$themeColors = $slide.ThemeColorScheme

Set-ThemeColor $themeColors 3  0x44 0x54 0x6A   # dk2
Set-ThemeColor $themeColors 4  0xE7 0xE6 0xE6   # lt2
Set-ThemeColor $themeColors 5  0x5B 0x9B 0xD5   # accent1
Set-ThemeColor $themeColors 6  0xED 0x7D 0x31   # accent2
Set-ThemeColor $themeColors 7  0xA5 0xA5 0xA5   # accent3
Set-ThemeColor $themeColors 8  0xFF 0xC0 0x00   # accent4
Set-ThemeColor $themeColors 9  0x44 0x72 0xC4   # accent5
Set-ThemeColor $themeColors 10 0x70 0xAD 0x47   # accent6
Set-ThemeColor $themeColors 11 0x05 0x63 0xC1   # hlink
Set-ThemeColor $themeColors 12 0x95 0x4F 0x72   # folHlink
